# Applies the diff:
#  1) Swaps the match-data columns (F:V) between 12 pairs of adjacent rows
#     (columns A-E - index/country/tournament/season/date - stay put).
#  2) Appends three brand-new match rows (171-173) at the bottom, copying
#     the formatting of the last existing data row (170) and then filling
#     in their values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows {
    param($sheet, [int]$row1, [int]$row2)

    $rng1 = $sheet.Range("F" + $row1 + ":V" + $row1)
    $rng2 = $sheet.Range("F" + $row2 + ":V" + $row2)

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value2 = $vals2
    $rng2.Value2 = $vals1
}

# --- 1) Swap the row pairs ---------------------------------------------
$pairs = @(
    @(17, 18),
    @(19, 20),
    @(22, 23),
    @(33, 34),
    @(42, 43),
    @(44, 45),
    @(84, 85),
    @(101, 102),
    @(134, 135),
    @(142, 143),
    @(156, 157),
    @(162, 163)
)

foreach ($pair in $pairs) {
    Swap-MatchRows $ws $pair[0] $pair[1]
}

# --- 2) Append new rows 171-173 -----------------------------------------
$newRows = @(
    @{ Row=171; A=170; E=45253.375;            F="Madura United"; G=1; H="Bali United";    I=2;
       J=2.05; K="21/11/2023 21:12"; L=2.56; M="23/11/2023 08:59";
       N=3.38; O="21/11/2023 21:12"; P=3.42; Q="23/11/2023 08:59";
       R=3.24; S="21/11/2023 21:12"; T=2.67; U="23/11/2023 08:59";
       V="https://www.betexplorer.com/football/indonesia/liga-1/madura-united-bali-united/AmoaKyep/" },

    @{ Row=172; A=171; E=45253.54166666666;    F="Persita";        G=3; H="RANS Nusantara"; I=0;
       J=2.09; K="22/11/2023 01:12"; L=2.59; M="23/11/2023 12:59";
       N=3.31; O="22/11/2023 01:12"; P=3.59; Q="23/11/2023 12:59";
       R=3.11; S="22/11/2023 01:12"; T=2.54; U="23/11/2023 12:59";
       V="https://www.betexplorer.com/football/indonesia/liga-1/persita-rans-nusantara/GnTQxpBB/" },

    @{ Row=173; A=172; E=45253.54166666666;    F="PSM Makassar";   G=3; H="Persikabo 1973"; I=1;
       J=1.44; K="22/11/2023 01:12"; L=1.38; M="23/11/2023 12:55";
       N=4.32; O="22/11/2023 01:12"; P=4.97; Q="23/11/2023 12:55";
       R=5.72; S="22/11/2023 01:12"; T=7.53; U="23/11/2023 12:55";
       V="https://www.betexplorer.com/football/indonesia/liga-1/psm-makassar-persikabo-1973/MPTMwQe5/" }
)

foreach ($nr in $newRows) {
    $r = $nr.Row

    # Bring over the formatting (bold/centered index col, date-time numeric
    # format on the match-date col, etc.) from the last existing data row.
    $ws.Range("A170:V170").Copy($ws.Range("A" + $r + ":V" + $r))

    $ws.Range("A$r").Value2 = $nr.A
    $ws.Range("B$r").Value2 = "indonesia"
    $ws.Range("C$r").Value2 = "liga-1"
    $ws.Range("D$r").Value2 = "2023-2024"
    $ws.Range("E$r").Value2 = $nr.E
    $ws.Range("F$r").Value2 = $nr.F
    $ws.Range("G$r").Value2 = $nr.G
    $ws.Range("H$r").Value2 = $nr.H
    $ws.Range("I$r").Value2 = $nr.I
    $ws.Range("J$r").Value2 = $nr.J
    $ws.Range("K$r").Value2 = $nr.K
    $ws.Range("L$r").Value2 = $nr.L
    $ws.Range("M$r").Value2 = $nr.M
    $ws.Range("N$r").Value2 = $nr.N
    $ws.Range("O$r").Value2 = $nr.O
    $ws.Range("P$r").Value2 = $nr.P
    $ws.Range("Q$r").Value2 = $nr.Q
    $ws.Range("R$r").Value2 = $nr.R
    $ws.Range("S$r").Value2 = $nr.S
    $ws.Range("T$r").Value2 = $nr.T
    $ws.Range("U$r").Value2 = $nr.U
    $ws.Range("V$r").Value2 = $nr.V
}
